$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the price/volume columns so Excel
# does not auto-convert numeric-looking strings (e.g. "13.60", "1.00",
# "89.873.44") into actual numbers, which would silently drop the
# formatting (trailing zeros / thousand-dot grouping) baked into the
# source text. The format is cleared again at the end so the saved
# cells end up with no explicit style, matching the original file.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range('D2').Value = '89.873.44'
$ws.Range('E2').Value = '  +3.29%  '
$ws.Range('D3').Value = '3.217.92'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '218.21'
$ws.Range('E5').Value = '  +5.62%  '
$ws.Range('D6').Value = '630.79'
$ws.Range('E6').Value = '  +4.20%  '
$ws.Range('D7').Value = '0.393'
$ws.Range('E7').Value = '  +8.04%  '
$ws.Range('D8').Value = '0.697'
$ws.Range('E8').Value = '  +6.65%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.213.79'
$ws.Range('E11').Value = '  +7.78%  '
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('E13').Value = '  +7.72%  '
$ws.Range('E14').Value = '  +3.65%  '
$ws.Range('E15').Value = '  +5.05%  '
$ws.Range('D16').Value = '3.811.87'
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').Value = '89.607.75'
$ws.Range('E17').Value = '  +3.28%  '
$ws.Range('D18').Value = '3.218.88'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '0.0000236'
$ws.Range('E19').Value = '  +83.61%  '
$ws.Range('D20').Value = '3.48'
$ws.Range('E20').Value = '  +18.24%  '
$ws.Range('D21').Value = '13.60'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').Value = '440.65'
$ws.Range('E22').Value = '  +7.29%  '
$ws.Range('D23').Value = '8.66'
$ws.Range('E23').Value = '  +2.60%  '
$ws.Range('D24').Value = '5.11'
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').Value = '12.04'
$ws.Range('E25').Value = '  +4.44%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = '5.26'
$ws.Range('E26').Value = '  +2.81%  '
$ws.Range('D27').Value = '82.75'
$ws.Range('E27').Value = '  +13.05%  '
$ws.Range('D28').Value = '3.391.45'
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('E32').Value = '  +39.98%  '
$ws.Range('E33').Value = '  +3.76%  '
$ws.Range('D34').Value = '546.41'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').Value = '7.06'
$ws.Range('E35').Value = '  +7.37%  '
$ws.Range('D36').Value = '1.93'
$ws.Range('E36').Value = '  +4.19%  '
$ws.Range('D37').Value = '1.31'
$ws.Range('E37').Value = '  +3.44%  '
$ws.Range('E38').Value = '  +3.76%  '
$ws.Range('E39').Value = '  +2.67%  '
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '1.94'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').Value = '146.14'
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('D46').Value = '174.00'
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('D47').Value = '43.78'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('D48').Value = '0.757'
$ws.Range('E48').Value = '  +10.18%  '
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('E51').Value = '  +6.85%  '

$fmtRange.ClearFormats()
